$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colombia Primera A")

# Row 14
$ws.Range("B14").Value = 6772177
$ws.Range("F14").Value = 'Aguilas Doradas'
$ws.Range("G14").Value = 'Alianza Petrolera'
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 'D'
$ws.Range("K14").Value = 2.15
$ws.Range("L14").Value = 3.3
$ws.Range("M14").Value = 3.5
$ws.Range("N14").Value = 2.2
$ws.Range("O14").Value = 3.5
$ws.Range("P14").Value = 3.2
$ws.Range("Q14").Value = -0.25
$ws.Range("R14").Value = 1.9
$ws.Range("S14").Value = 1.9
$ws.Range("T14").Value = 2.75
$ws.Range("U14").Value = 1.95
$ws.Range("V14").Value = 1.85
$ws.Range("W14").Value = -1
$ws.Range("X14").Value = 2.5
$ws.Range("Z14").Value = -0.5
$ws.Range("AA14").Value = 0.45
$ws.Range("AB14").Value = -1
$ws.Range("AC14").Value = 0.8500000000000001

# Row 15
$ws.Range("B15").Value = 6772175
$ws.Range("F15").Value = 'Atletico Nacional Medellin'
$ws.Range("G15").Value = 'Deportivo Pasto'
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 'H'
$ws.Range("K15").Value = 1.666
$ws.Range("L15").Value = 3.75
$ws.Range("M15").Value = 4.5
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 3.6
$ws.Range("P15").Value = 5
$ws.Range("Q15").Value = -0.75
$ws.Range("R15").Value = 2
$ws.Range("S15").Value = 1.85
$ws.Range("T15").Value = 2.25
$ws.Range("U15").Value = 1.85
$ws.Range("V15").Value = 2
$ws.Range("W15").Value = 0.8
$ws.Range("X15").Value = -1
$ws.Range("Z15").Value = 0.5
$ws.Range("AA15").Value = -0.5
$ws.Range("AB15").Value = 0.8500000000000001
$ws.Range("AC15").Value = -1

# Row 208
$ws.Range("B208").Value = 7404214
$ws.Range("F208").Value = 'Boyaca Chico'
$ws.Range("G208").Value = 'Deportivo Cali'
$ws.Range("H208").Value = 1
$ws.Range("J208").Value = 'D'
$ws.Range("K208").Value = 3.2
$ws.Range("M208").Value = 2.2
$ws.Range("N208").Value = 3.6
$ws.Range("O208").Value = 3
$ws.Range("P208").Value = 2.25
$ws.Range("Q208").Value = 0.25
$ws.Range("R208").Value = 1.95
$ws.Range("T208").Value = 2.25
$ws.Range("U208").Value = 1.875
$ws.Range("V208").Value = 1.975
$ws.Range("X208").Value = 2
$ws.Range("Y208").Value = -1
$ws.Range("Z208").Value = 0.475
$ws.Range("AA208").Value = -0.5
$ws.Range("AB208").Value = -0.5
$ws.Range("AC208").Value = 0.4875

# Row 211
$ws.Range("B211").Value = 7404217
$ws.Range("F211").Value = 'Alianza Petrolera'
$ws.Range("G211").Value = 'Deportivo Pereira'
$ws.Range("H211").Value = 2
$ws.Range("J211").Value = 'H'
$ws.Range("K211").Value = 1.95
$ws.Range("L211").Value = 3.2
$ws.Range("M211").Value = 3.75
$ws.Range("N211").Value = 1.95
$ws.Range("O211").Value = 3.2
$ws.Range("P211").Value = 4.75
$ws.Range("Q211").Value = -0.5
$ws.Range("R211").Value = 1.925
$ws.Range("S211").Value = 1.875
$ws.Range("T211").Value = 2
$ws.Range("U211").Value = 1.825
$ws.Range("W211").Value = 0.95
$ws.Range("X211").Value = -1
$ws.Range("Z211").Value = 0.925
$ws.Range("AA211").Value = -1
$ws.Range("AB211").Value = 0.825
$ws.Range("AC211").Value = -1

# Row 212
$ws.Range("B212").Value = 7404216
$ws.Range("F212").Value = 'Independiente Santa Fe'
$ws.Range("G212").Value = 'Once Caldas'
$ws.Range("H212").Value = 0
$ws.Range("J212").Value = 'A'
$ws.Range("K212").Value = 1.85
$ws.Range("L212").Value = 3.1
$ws.Range("M212").Value = 4.2
$ws.Range("N212").Value = 2.25
$ws.Range("P212").Value = 3.3
$ws.Range("Q212").Value = -0.25
$ws.Range("R212").Value = 1.9
$ws.Range("S212").Value = 1.9
$ws.Range("T212").Value = 2.5
$ws.Range("U212").Value = 1.925
$ws.Range("V212").Value = 1.925
$ws.Range("W212").Value = -1
$ws.Range("Y212").Value = 2.3
$ws.Range("Z212").Value = -1
$ws.Range("AA212").Value = 0.8999999999999999
$ws.Range("AB212").Value = -1
$ws.Range("AC212").Value = 0.925

# Row 214
$ws.Range("B214").Value = 7404215
$ws.Range("F214").Value = 'America de Cali'
$ws.Range("G214").Value = 'Atletico Bucaramanga'
$ws.Range("H214").Value = 1
$ws.Range("I214").Value = 2
$ws.Range("K214").Value = 1.444
$ws.Range("L214").Value = 4.5
$ws.Range("M214").Value = 6
$ws.Range("N214").Value = 1.363
$ws.Range("O214").Value = 5
$ws.Range("P214").Value = 7.5
$ws.Range("Q214").Value = -1.25
$ws.Range("R214").Value = 1.775
$ws.Range("S214").Value = 2.025
$ws.Range("T214").Value = 3
$ws.Range("U214").Value = 1.925
$ws.Range("V214").Value = 1.875
$ws.Range("Y214").Value = 6.5
$ws.Range("AA214").Value = 1.025
$ws.Range("AB214").Value = 0
$ws.Range("AC214").Value = -0

# Row 216
$ws.Range("B216").Value = 7404260
$ws.Range("F216").Value = 'Atletico Nacional Medellin'
$ws.Range("G216").Value = 'Deportes Tolima'
$ws.Range("H216").Value = 2
$ws.Range("I216").Value = 3
$ws.Range("K216").Value = 2
$ws.Range("L216").Value = 3.25
$ws.Range("M216").Value = 3.5
$ws.Range("N216").Value = 1.75
$ws.Range("O216").Value = 3.6
$ws.Range("P216").Value = 4.75
$ws.Range("Q216").Value = -0.75
$ws.Range("R216").Value = 2
$ws.Range("S216").Value = 1.8
$ws.Range("T216").Value = 2.5
$ws.Range("U216").Value = 2
$ws.Range("V216").Value = 1.8
$ws.Range("Y216").Value = 3.75
$ws.Range("AA216").Value = 0.8
$ws.Range("AB216").Value = 1
$ws.Range("AC216").Value = -1

# Row 217
$ws.Range("B217").Value = 7404213
$ws.Range("F217").Value = 'Jaguares de Cordoba'
$ws.Range("G217").Value = 'Aguilas Doradas'
$ws.Range("H217").Value = 0
$ws.Range("I217").Value = 1
$ws.Range("K217").Value = 3.25
$ws.Range("L217").Value = 3.1
$ws.Range("M217").Value = 2.2
$ws.Range("N217").Value = 3.6
$ws.Range("O217").Value = 3.2
$ws.Range("P217").Value = 2.15
$ws.Range("Q217").Value = 0.25
$ws.Range("R217").Value = 1.975
$ws.Range("S217").Value = 1.825
$ws.Range("T217").Value = 2
$ws.Range("U217").Value = 1.75
$ws.Range("V217").Value = 2.05
$ws.Range("Y217").Value = 1.15
$ws.Range("AA217").Value = 0.825
$ws.Range("AB217").Value = -1
$ws.Range("AC217").Value = 1.05

# Row 238
$ws.Range("B238").Value = 7528604
$ws.Range("F238").Value = 'Aguilas Doradas'
$ws.Range("G238").Value = 'Deportivo Cali'
$ws.Range("H238").Value = 3
$ws.Range("J238").Value = 'H'
$ws.Range("K238").Value = 1.666
$ws.Range("L238").Value = 3.75
$ws.Range("M238").Value = 5
$ws.Range("N238").Value = 1.363
$ws.Range("O238").Value = 5
$ws.Range("P238").Value = 9
$ws.Range("Q238").Value = -1.25
$ws.Range("R238").Value = 1.825
$ws.Range("T238").Value = 2.75
$ws.Range("U238").Value = 1.9
$ws.Range("V238").Value = 1.9
$ws.Range("W238").Value = 0.363
$ws.Range("Y238").Value = -1
$ws.Range("Z238").Value = 0.825
$ws.Range("AA238").Value = -1
$ws.Range("AB238").Value = 0.8999999999999999
$ws.Range("AC238").Value = -1

# Row 239
$ws.Range("B239").Value = 7528136
$ws.Range("F239").Value = 'Millonarios'
$ws.Range("G239").Value = 'Atletico Nacional Medellin'
$ws.Range("H239").Value = 0
$ws.Range("J239").Value = 'A'
$ws.Range("K239").Value = 1.85
$ws.Range("L239").Value = 3.3
$ws.Range("M239").Value = 4.5
$ws.Range("N239").Value = 1.85
$ws.Range("O239").Value = 3.5
$ws.Range("P239").Value = 4.2
$ws.Range("Q239").Value = -0.5
$ws.Range("R239").Value = 1.875
$ws.Range("T239").Value = 2.5
$ws.Range("U239").Value = 2.05
$ws.Range("V239").Value = 1.8
$ws.Range("W239").Value = -1
$ws.Range("Y239").Value = 3.2
$ws.Range("Z239").Value = -1
$ws.Range("AA239").Value = 0.9750000000000001
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.8

# Row 392
$ws.Range("R392").Value = 2
$ws.Range("S392").Value = 1.85
$ws.Range("U392").Value = 1.95
$ws.Range("V392").Value = 1.9

# Row 393
$ws.Range("N393").Value = 3.1
$ws.Range("P393").Value = 2.3
$ws.Range("R393").Value = 1.825
$ws.Range("S393").Value = 2.025
$ws.Range("T393").Value = 2.5
$ws.Range("U393").Value = 2.05
$ws.Range("V393").Value = 1.8

# Row 394
$ws.Range("R394").Value = 1.975
$ws.Range("S394").Value = 1.875

# Row 395
$ws.Range("R395").Value = 1.85
$ws.Range("S395").Value = 2

# Row 396
$ws.Range("R396").Value = 1.9
$ws.Range("S396").Value = 1.95

# Row 398
$ws.Range("N398").Value = 1.95
$ws.Range("P398").Value = 4.2
$ws.Range("R398").Value = 2.025
$ws.Range("S398").Value = 1.825
$ws.Range("U398").Value = 1.875
$ws.Range("V398").Value = 1.975

# Row 400
$ws.Range("R400").Value = 1.85
$ws.Range("S400").Value = 2
$ws.Range("U400").Value = 1.975
$ws.Range("V400").Value = 1.875

# Row 401
$ws.Range("U401").Value = 1.8
$ws.Range("V401").Value = 2.05
